$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.758.13'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.339.23'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.26'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.02%  '
$ws.Range('E7').Value = '  +1.12%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.620'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.12'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.62'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.01'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.108'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.52'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.691.07'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.341.81'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.661.72'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.61'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.14'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.64%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.26'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.59'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.69'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +7.32%  '
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.88'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.72'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '167.96'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0887'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.77'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +7.06%  '
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.76'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.41%  '
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  +4.35%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.83'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.87'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.29%  '
$ws.Range('E41').Value = '  +9.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '104.64'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +13.41%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.237'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.11%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.43'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.02%  '
$ws.Range('B45').Value = 'MultiversX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '71.74'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.07%  '
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '114.26'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.664.10'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '77.24'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('B50').Value = 'MinaProtocolToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.57'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +11.56%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.97'
$ws.Range('D51').ClearFormats()
